# Adds "onsite" / "remote" / "blended" delivery-mode columns (O, P, Q) to the
# programlisting sheet, filling every program row with 1/1/1 (available in all
# modes) except the last row (IT Tech Support Specialist Diploma), which is
# onsite-only in the source data (0/1/0). This also overwrites the stray
# "move *csw_with_aw*..." note that used to live in O18, which naturally
# drops that now-unused string out of the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column headers
$ws.Range("O1").Value = "onsite"
$ws.Range("P1").Value = "remote"
$ws.Range("Q1").Value = "blended"

# Data rows 2-51: available onsite, remote, and blended
$ws.Range("O2:Q51").Value = 1

# Last data row (52): onsite only
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 1
$ws.Range("Q52").Value = 0

# Reflect the new used range in the view: scroll near the bottom-right and
# select the next empty cell, as a user would after finishing the edit.
$ws.Range("R52").Select()
